$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D2:D6 to "U"
$ws.Range("D2:D6").Value = "U"

# Clear F2:G6 contents (they are no longer used)
$ws.Range("F2:G6").ClearContents()

# Set selection to D7 to match the saved view state
$ws.Range("D7").Select()
